# Updates the cryptos price table (columns D = Price, E = Volume(1h))
# as scraped by the GitHub Actions job. Numeric-looking Price strings are
# prefixed with a leading apostrophe so Excel stores them as text, just
# like the original sheet (e.g. "1.000" must stay literal text, not 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.207.08"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "'1.852.58"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'0.6990"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").Value = "'236.99"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D8").Value = "'0.07883"
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").Value = "'23.90"
$ws.Range("E10").Value = "  +3.17%  "
$ws.Range("D11").Value = "'0.08133"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "'1.849.24"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "'5.184"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'0.7052"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").Value = "'89.32"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "'29.213.10"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'5.797"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "'0.000007821"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "'235.25"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'2.094.84"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'7.486"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "'162.29"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").Value = "'8.885"
$ws.Range("D27").Value = "'0.1416"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "'18.01"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'1.906"
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").Value = "'1.476"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "'4.292"
$ws.Range("E32").Value = "  -4.68%  "
$ws.Range("D33").Value = "'4.009"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "'0.05148"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "'0.7057"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").Value = "'0.9976"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'2.704"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").Value = "'1.149.16"
$ws.Range("E41").Value = "  +4.16%  "
$ws.Range("D42").Value = "'0.9204"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("D43").Value = "'5.939"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").Value = "'0.4232"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "'69.89"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'102.83"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("D48").Value = "'0.5293"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("D49").Value = "'1.736"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").Value = "'9.149"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'6.948"
$ws.Range("E51").Value = "  -0.59%  "
